# Update cryptos list (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "71.132.89"
$ws.Range("E2").Value = "  -1.71%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.569.27"
$ws.Range("E3").Value = "  -5.43%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
$ws.Range("D5").Value = "581.98"
$ws.Range("E5").Value = "  -3.22%  "

# Row 6 - Solana
$ws.Range("E6").Value = "  -3.23%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -2.50%  "

# Row 9 - becomes Dogecoin (was LidoStakedEther)
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.167"
$ws.Range("E9").Value = "  -1.56%  "

# Row 10 - becomes LidoStakedEther (was Dogecoin)
$ws.Range("B10").Value = "LidoStakedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D10").Value = "2.570.14"
$ws.Range("E10").Value = "  -5.36%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.21%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "0.351"
$ws.Range("E12").Value = "  -1.21%  "

# Row 13 - Toncoin
$ws.Range("D13").Value = "4.85"
$ws.Range("E13").Value = "  -3.53%  "

# Row 14 - becomes ShibaInu (was WrappedliquidstakedEther2.0)
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "0.0000186"
$ws.Range("E14").Value = "  -0.45%  "

# Row 15 - becomes WrappedliquidstakedEther2.0 (was ShibaInu)
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.054.92"
$ws.Range("E15").Value = "  -4.81%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "70.945.18"
$ws.Range("E16").Value = "  -1.67%  "

# Row 17 - Avalanche
$ws.Range("D17").Value = "25.28"
$ws.Range("E17").Value = "  -4.36%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.579.03"
$ws.Range("E18").Value = "  -4.91%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "11.91"
$ws.Range("E19").Value = "  -3.64%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "7.69"
$ws.Range("E20").Value = "  -5.89%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "365.05"
$ws.Range("E21").Value = "  -2.66%  "

# Row 22 - Polkadot
$ws.Range("E22").Value = "  -4.07%  "

# Row 23 - SuiNetwork
$ws.Range("D23").Value = "2.00"
$ws.Range("E23").Value = "  -1.83%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  -0.08%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "70.47"
$ws.Range("E25").Value = "  -2.74%  "

# Row 26 - NEARProtocol
$ws.Range("E26").Value = "  -5.75%  "

# Row 27 - Aptos
$ws.Range("D27").Value = "9.27"
$ws.Range("E27").Value = "  -6.46%  "

# Row 28 - WrappedeETH
$ws.Range("D28").Value = "2.755.24"
$ws.Range("E28").Value = "  -3.48%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  +0.70%  "

# Row 30 - PEPE
$ws.Range("E30").Value = "  -7.57%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "7.82"
$ws.Range("E31").Value = "  -4.40%  "

# Row 32 - becomes Bittensor (was Fetch.AI)
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "486.95"
$ws.Range("E32").Value = "  -4.68%  "

# Row 33 - becomes Fetch.AI (was Bittensor)
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "1.31"
$ws.Range("E33").Value = "  -1.29%  "

# Row 34 - PancakeSwap
$ws.Range("E34").Value = "  -2.85%  "

# Row 35 - FirstDigitalUSD
$ws.Range("E35").Value = "  +0.11%  "

# Row 36 - Monero
$ws.Range("D36").Value = "156.23"
$ws.Range("E36").Value = "  -4.88%  "

# Row 37 - Kaspa
$ws.Range("E37").Value = "  +4.64%  "

# Row 38 - EthereumClassic
$ws.Range("D38").Value = "18.89"
$ws.Range("E38").Value = "  -4.47%  "

# Row 39 - WhiteBITCoin
$ws.Range("D39").Value = "18.86"
$ws.Range("E39").Value = "  -1.33%  "

# Row 40 - ImmutableX
$ws.Range("E40").Value = "  -4.47%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  -6.33%  "

# Row 43 - dogwifhat
$ws.Range("D43").Value = "2.48"
$ws.Range("E43").Value = "  -4.07%  "

# Row 44 - RenderToken
$ws.Range("E44").Value = "  -6.28%  "

# Row 45 - PolygonEcosystemToken
$ws.Range("D45").Value = "0.321"
$ws.Range("E45").Value = "  -4.67%  "

# Row 46 - OKB
$ws.Range("D46").Value = "38.64"
$ws.Range("E46").Value = "  -2.25%  "

# Row 47 - Aave
$ws.Range("D47").Value = "147.38"
$ws.Range("E47").Value = "  -6.25%  "

# Row 48 - Filecoin
$ws.Range("E48").Value = "  -4.82%  "

# Row 49 - ARBITRUM
$ws.Range("D49").Value = "0.533"
$ws.Range("E49").Value = "  -6.30%  "

# Row 50 - Optimism
$ws.Range("E50").Value = "  -8.32%  "

# Row 51 - Mantle
$ws.Range("D51").Value = "0.596"
$ws.Range("E51").Value = "  -2.39%  "
